$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("客户档案")

# Update data row cells first so new shared strings are appended in the
# same order as the target workbook ("北中国" before "南北中国").
$ws.Range("A2").Value = "北中国"

# Update header cell A1 ("团队" -> "南北中国")
$ws.Range("A1").Value = "南北中国"

$ws.Range("D2").Value = "测试1"

# Update the active selection shown in the sheet view (was K24, now D3)
$ws.Range("D3").Select()
